# Insert a new data row at row 33 (pushes existing rows 33-49 down to 34-50)
# and populate it with the new record's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

$ws.Range("A33").Value() = 10
$ws.Range("B33").Value() = 'Vega Modelo de Temuco'
$ws.Range("C33").Value() = 'La Araucanía'
$ws.Range("D33").Value() = 44762
$ws.Range("E33").Value() = 9
$ws.Range("F33").Value() = 'Fruta'
$ws.Range("G33").Value() = 100108
$ws.Range("H33").Value() = 'Tropicales y subtropicales'
$ws.Range("I33").Value() = 100108003
$ws.Range("J33").Value() = 'Maracuyá'
$ws.Range("K33").Value() = 'Sin especificar'
$ws.Range("L33").Value() = 'Primera'
$ws.Range("M33").Value() = 50
$ws.Range("N33").Value() = 34000
$ws.Range("O33").Value() = 34000
$ws.Range("P33").Value() = 34000
$ws.Range("Q33").Value() = '$/caja 18 kilos'
$ws.Range("R33").Value() = 'Región de Arica y Parinacota'
$ws.Range("S33").Value() = 1889
$ws.Range("T33").Value() = 18
